$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
$th = $nm.Theme
$te = $th.ThemeElements
$fs = $te.ThemeFontScheme
$major = $fs.MajorFont
$major.Latin = "Verdana"
Write-Host "set done"
